$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.153.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.95%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.413.61'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +7.25%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '293.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.72%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.89%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.557'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.38%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.09%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.499'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.68%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.66%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0774'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.15%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.00'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.63%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.51%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.767.42'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.78%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.389.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.11%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.08%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.831'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.43%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.126.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.99%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.33%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0937'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.03%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.52%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.63'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.45%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.52%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.23%  '

# Row 25
$ws.Range('E25').Value = '  +0.09%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.43%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.69%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.44'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -11.15%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.95%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.81'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +21.48%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '21.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.38%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.21%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.41%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.39'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.34%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0758'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.88%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +17.66%  '

# Row 37
$ws.Range('E37').Value = '  -1.71%  '

# Row 38
$ws.Range('E38').Value = '  -0.79%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.75%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.72'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.96%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0291'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.36%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.984.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.67%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.74%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.04%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '89.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.81%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -12.30%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +23.51%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.37'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.46%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '99.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.10%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.637.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.84%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.182'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.79%  '

